# Insert a new row above the current row 13 ("pointing offset angle S/C")
# for a new "S/C-Sun distance" parameter, pushing all subsequent rows
# (and the trailing blank formatting row) down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(13).Insert() | Out-Null

# New row 13 content: S/C-Sun distance (km). Not applicable for the
# Moon / Mars missions (columns C/D) or the overall "required" column H;
# Venus / Europa / Earth (columns E/F/G) get the actual mean distances.
$ws.Range("A13").Value = "S/C-Sun distance"
$ws.Range("B13").Value = "km"
$ws.Range("C13").Value = "N/A"
$ws.Range("D13").Value = "N/A"
$ws.Range("E13").Value = 227900000
$ws.Range("F13").Value = 108200000
$ws.Range("G13").Value = 780000000
$ws.Range("H13").Value = "N/A"

# Distances get the highlighted-fill scientific number format.
$ws.Range("E13:G13").Interior.Color = $ws.Range("C12").Interior.Color
$ws.Range("E13:G13").NumberFormat = "0.00E+00"

# Match the author's final selection on the newly inserted row.
$ws.Range("D13").Select() | Out-Null
